$wb = $excel.ActiveWorkbook

# --- "Lop" sheet: add a new column D (Khoa abbreviation "CT") ---
$wsLop = $wb.Worksheets.Item("Lop")
$wsLop.Range("D1").Value = "CT"
$wsLop.Range("D2").Value = "CT"
$wsLop.Range("D3").Value = "CT"

# --- "SinhVien" sheet: column A becomes the numeric student id,   ---
# --- the class name that used to live in A now moves to column J ---
$wsSV = $wb.Worksheets.Item("SinhVien")
$wsSV.Range("J1").Value = $wsSV.Range("A1").Value2
$wsSV.Range("J2").Value = $wsSV.Range("A2").Value2

$wsSV.Range("A1").Value = 118001525
$wsSV.Range("A2").Value = 118001526
$wsSV.Range("A3").Value = 118001527

# Column width tweaks on SinhVien (D, E, F got narrower)
$wsSV.Columns.Item(4).ColumnWidth = 23.33
$wsSV.Columns.Item(5).ColumnWidth = 9.95
$wsSV.Columns.Item(6).ColumnWidth = 10.45

# Selections / active sheet: "Lop" becomes the active tab with D3
# selected, "SinhVien" is left with A4 selected.
[void]$wsSV.Range("A4").Select()
$wsLop.Select()
[void]$wsLop.Range("D3").Select()
